$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.879.55"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3
$ws.Range("D3").Value = "'3.523.50"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'602.03"
$ws.Range("E5").Value = "  -1.60%  "

# Row 6
$ws.Range("D6").Value = "'195.55"
$ws.Range("E6").Value = "  +5.51%  "

# Row 7
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("E9").Value = "  -2.27%  "

# Row 10
$ws.Range("D10").Value = "'0.654"
$ws.Range("E10").Value = "  +0.85%  "

# Row 11
$ws.Range("D11").Value = "'54.03"
$ws.Range("E11").Value = "  +0.71%  "

# Row 12
$ws.Range("D12").Value = "'0.0000302"
$ws.Range("E12").Value = "  -2.72%  "

# Row 13
$ws.Range("D13").Value = "'9.53"
$ws.Range("E13").Value = "  +0.64%  "

# Row 14
$ws.Range("D14").Value = "'4.077.14"
$ws.Range("E14").Value = "  -0.65%  "

# Row 15
$ws.Range("D15").Value = "'602.88"
$ws.Range("E15").Value = "  -3.29%  "

# Row 16
$ws.Range("D16").Value = "'70.039.79"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("D17").Value = "'19.13"
$ws.Range("E17").Value = "  +1.32%  "

# Row 18
$ws.Range("D18").Value = "'12.63"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19
$ws.Range("D19").Value = "'3.517.14"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("E20").Value = "  +0.64%  "

# Row 21
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
$ws.Range("D22").Value = "'18.31"
$ws.Range("E22").Value = "  +3.63%  "

# Row 23
$ws.Range("D23").Value = "'5.24"
$ws.Range("E23").Value = "  +6.70%  "

# Row 24
$ws.Range("D24").Value = "'103.72"
$ws.Range("E24").Value = "  +2.41%  "

# Row 25
$ws.Range("D25").Value = "'4.61"
$ws.Range("E25").Value = "  -2.55%  "

# Row 26
$ws.Range("E26").Value = "  +2.44%  "

# Row 27
$ws.Range("D27").Value = "'10.93"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
$ws.Range("D28").Value = "'9.68"
$ws.Range("E28").Value = "  +1.02%  "

# Row 29
$ws.Range("D29").Value = "'33.63"
$ws.Range("E29").Value = "  +0.23%  "

# Row 30
$ws.Range("D30").Value = "'4.50"
$ws.Range("E30").Value = "  +24.90%  "

# Row 31
$ws.Range("D31").Value = "'7.12"
$ws.Range("E31").Value = "  +1.17%  "

# Row 32
$ws.Range("D32").Value = "'12.71"
$ws.Range("E32").Value = "  +3.82%  "

# Row 33
$ws.Range("E33").Value = "  +1.62%  "

# Row 34
$ws.Range("D34").Value = "'63.20"
$ws.Range("E34").Value = "  -0.45%  "

# Row 35
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "'0.0₃0828"
$ws.Range("E35").Value = "  +5.98%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "'3.766.54"
$ws.Range("E36").Value = "  +6.33%  "

# Row 37
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'3.08"
$ws.Range("E37").Value = "  -4.98%  "

# Row 38
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.07%  "

# Row 39
$ws.Range("D39").Value = "'0.393"
$ws.Range("E39").Value = "  -1.87%  "

# Row 40
$ws.Range("D40").Value = "'3.59"
$ws.Range("E40").Value = "  +0.86%  "

# Row 41
$ws.Range("D41").Value = "'36.76"
$ws.Range("E41").Value = "  -1.03%  "

# Row 42
$ws.Range("D42").Value = "'489.29"
$ws.Range("E42").Value = "  -8.05%  "

# Row 43
$ws.Range("E43").Value = "  -0.44%  "

# Row 44
$ws.Range("E44").Value = "  -0.37%  "

# Row 45
$ws.Range("E45").Value = "  -1.98%  "

# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.33"
$ws.Range("E46").Value = "  -1.03%  "

# Row 47
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.83"
$ws.Range("E47").Value = "  -3.61%  "

# Row 48
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("D49").Value = "'8.66"
$ws.Range("E49").Value = "  -5.61%  "

# Row 50
$ws.Range("D50").Value = "'1.34"
$ws.Range("E50").Value = "  +14.01%  "

# Row 51
$ws.Range("D51").Value = "'0.000244"
$ws.Range("E51").Value = "  +1.13%  "
